# Update the COVID-19 Valais daily figures table on sheet "Feuil1".
# Rows 355-357 receive updated / newly-entered daily figures. Columns
# B (Cumul cas positifs), H (Total hospitalisations), J (Cumul deces) and
# K (Nb nouveaux deces) are driven by worksheet formulas already present in
# the sheet, so they recompute automatically once the raw inputs below are
# written - we only need to touch the source columns (C, E, F, G, L, M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 355 (2020-12-02): revised case count + one new extra-hospital death
$ws.Cells.Item(355, 3).Value = 33   # C355 Nb nouveaux cas positifs: 31 -> 33

# Column M (Nb nouveaux deces extra-hospitaliers) is formatted as Text
# (numFmtId 49, "@"). Writing a plain number through .Value on a Text
# formatted cell stores it as text, so flip the cell to a general number
# format for the write and restore the original text format afterwards -
# this mirrors how the source values in this column are genuine numbers.
$ws.Cells.Item(355, 13).NumberFormat = "General"
$ws.Cells.Item(355, 13).Value = 1   # M355 Nb nouveaux deces extra-hosp.: 0 -> 1
$ws.Cells.Item(355, 13).NumberFormat = "@"

# --- Row 356 (2020-12-03): revised case count
$ws.Cells.Item(356, 3).Value = 46   # C356 Nb nouveaux cas positifs: 10 -> 46

# --- Row 357 (2020-12-04): first entry of the day's figures (was blank)
$ws.Cells.Item(357, 3).Value = 7    # C357 Nb nouveaux cas positifs
$ws.Cells.Item(357, 5).Value = 10   # E357 Patients COVID-19 aux SI total
$ws.Cells.Item(357, 6).Value = 7    # F357 Patients COVID-19 intubes
$ws.Cells.Item(357, 7).Value = 41   # G357 Patients hospitalises hors SI

$ws.Cells.Item(357, 12).NumberFormat = "General"
$ws.Cells.Item(357, 12).Value = 0   # L357 Nb nouveaux deces a l'hopital
$ws.Cells.Item(357, 12).NumberFormat = "@"

$ws.Cells.Item(357, 13).NumberFormat = "General"
$ws.Cells.Item(357, 13).Value = 0   # M357 Nb nouveaux deces extra-hosp.
$ws.Cells.Item(357, 13).NumberFormat = "@"
